$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "servicios"

$headerRange = $ws.Range("A1:E1")
$headerRange.ClearFormats()

$cols = $ws.Range("A1:E1").EntireColumn
$cols.ColumnWidth = 8
